# DG: update model API
#
# 1) Refresh the "last saved" date shown by the datetimeFigureOut fields on
#    the slide master, every slide layout, and the notes master.
# 2) Rename the model call from deletePerson(p) to deletePersons(p) in the
#    sequence diagram on slide 1 (run split moves after "deletePersons(p").

$p = $ppt.ActivePresentation
$oldDate = "10/16/2016"
$newDate = "06-Jan-17"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            $tr = $sh.TextFrame.TextRange
            $found = $tr.Find($oldDate, 0)
            if ($found -ne $null) {
                $found.Text = $newDate
            } else {
                $tr.Text = $newDate
            }
        }
    }
}

# Slide master
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every slide layout hanging off the master
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# Notes master
$notesMaster = $p.NotesMaster
Update-DatePlaceholder $notesMaster.Shapes

# Slide 1: deletePerson(p) -> deletePersons(p)
$slide = $p.Slides.Item(1)
$callout = $slide.Shapes.Item("TextBox 28")
$tr = $callout.TextFrame.TextRange

$tail = $tr.Find("(p)", 0)
$tail.Text = ")"

$head = $slide.Shapes.Item("TextBox 28").TextFrame.TextRange.Find("deletePerson", 0)
$head.Text = "deletePersons(p"
